$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for new columns
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

# Fill I (always 1) and J (copy of H) for rows 2 through 36
for ($r = 2; $r -le 36; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
